$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the existing "Eventos" long description (row 54, column C)
$ws.Range('C54').Value = 'Son los eventos deportivos que aparecerán en la aplicaciión. Se almacenará: Nombre, Deporte, Lugar, Organizador, participantes, pagina oficial, resultados y realizado (sí o no).'

# 2. Row 47 grows taller (wraps to one more line) after the edit
$ws.Rows(47).RowHeight = 37.45

# 3. Append new requirement rows 61-104
# Row 61
$ws.Range('A61').Value = 'R60'
$ws.Range('B61').Value = 'Recuperar contraseña'
$ws.Range('C61').Value = 'Un usuario puede cambiar su contraseña si no la recuerda'
$ws.Range('D61').Value = 'Importante'
$ws.Range('E61').Value = 'Funcional'
$ws.Range('F61').Value = 'Media'
$ws.Range('G61').Value = 'v2'
$ws.Rows(61).RowHeight = 12.8

# Row 62
$ws.Range('A62').Value = 'R61'
$ws.Range('B62').Value = 'Buscar eventos por nombre'
$ws.Range('C62').Value = 'Se podrá buscar eventos por su nombre'
$ws.Range('D62').Value = 'Importante'
$ws.Range('E62').Value = 'Funcional'
$ws.Range('F62').Value = 'Fácil'
$ws.Range('G62').Value = 'v2'
$ws.Rows(62).RowHeight = 23.85

# Row 63
$ws.Range('A63').Value = 'R62'
$ws.Range('B63').Value = 'Buscar eventos por tipo de deporte'
$ws.Range('C63').Value = 'Se buscará eventos por el tipo de deporte que se indique'
$ws.Range('D63').Value = 'Importante'
$ws.Range('E63').Value = 'Funcional'
$ws.Range('F63').Value = 'Fácil'
$ws.Range('G63').Value = 'v2'
$ws.Rows(63).RowHeight = 23.85

# Row 64
$ws.Range('A64').Value = 'R63'
$ws.Range('B64').Value = 'Chat'
$ws.Range('C64').Value = 'La web mostrará un chat general para todos los usuarios logueados'
$ws.Range('D64').Value = 'Opcional'
$ws.Range('E64').Value = 'Funcional'
$ws.Range('F64').Value = 'Difícil'
$ws.Range('G64').Value = 'v3'
$ws.Rows(64).RowHeight = 23.85

# Row 65
$ws.Range('A65').Value = 'R64'
$ws.Range('B65').Value = 'Salas de chat'
$ws.Range('C65').Value = 'Cada evento dispondrá de una sala de chat.'
$ws.Range('D65').Value = 'Opcional'
$ws.Range('E65').Value = 'Funcional'
$ws.Range('F65').Value = 'Difícil'
$ws.Range('G65').Value = 'v3'
$ws.Rows(65).RowHeight = 12.8

# Row 66
$ws.Range('A66').Value = 'R65'
$ws.Range('B66').Value = 'Valorar evento'
$ws.Range('C66').Value = 'Se podrá valorar el evento con un sistema de votación que irá del 1 al 5.'
$ws.Range('D66').Value = 'Importante'
$ws.Range('E66').Value = 'Funcional'
$ws.Range('F66').Value = 'Media'
$ws.Range('G66').Value = 'v2'
$ws.Rows(66).RowHeight = 23.85

# Row 67
$ws.Range('A67').Value = 'R66'
$ws.Range('B67').Value = 'Ranking de usuarios creadores'
$ws.Range('C67').Value = 'Ranking donde se mostrará los usuarios con más eventos creados'
$ws.Range('D67').Value = 'Opcional'
$ws.Range('E67').Value = 'Funcional'
$ws.Range('F67').Value = 'Fácil'
$ws.Range('G67').Value = 'v2'
$ws.Rows(67).RowHeight = 23.85

# Row 68
$ws.Range('A68').Value = 'R67'
$ws.Range('B68').Value = 'Ranking de usuarios con más comentarios'
$ws.Range('C68').Value = 'Ranking donde se mostrará los usuarios con más comentarios'
$ws.Range('D68').Value = 'Opcional'
$ws.Range('E68').Value = 'Funcional'
$ws.Range('F68').Value = 'Fácil'
$ws.Range('G68').Value = 'v2'
$ws.Rows(68).RowHeight = 23.85

# Row 69
$ws.Range('A69').Value = 'R68'
$ws.Range('B69').Value = 'Ranking de usuarios con más votos'
$ws.Range('C69').Value = 'Ranking donde se mostrará los usuarios con más votos a los eventos'
$ws.Range('D69').Value = 'Opcional'
$ws.Range('E69').Value = 'Funcional'
$ws.Range('F69').Value = 'Fácil'
$ws.Range('G69').Value = 'v2'
$ws.Rows(69).RowHeight = 23.85

# Row 70
$ws.Range('A70').Value = 'R69'
$ws.Range('B70').Value = 'Recibir notificaciones'
$ws.Range('C70').Value = 'Si el usuario lo desea, se mandará un correo cada vez que haya un nuevo evento.'
$ws.Range('D70').Value = 'Opcional'
$ws.Range('E70').Value = 'Funcional'
$ws.Range('F70').Value = 'Difícil'
$ws.Range('G70').Value = 'v3'
$ws.Rows(70).RowHeight = 23.85

# Row 71
$ws.Range('A71').Value = 'R70'
$ws.Range('B71').Value = 'Seguir deporte'
$ws.Range('C71').Value = 'Se podrá seguir un deporte y recibir notificaciones de los eventos de ese deporte.'
$ws.Range('D71').Value = 'Opcional'
$ws.Range('E71').Value = 'Funcional'
$ws.Range('F71').Value = 'Difícil'
$ws.Range('G71').Value = 'v3'
$ws.Rows(71).RowHeight = 23.85

# Row 72
$ws.Range('A72').Value = 'R71'
$ws.Range('B72').Value = 'Deporte'
$ws.Range('C72').Value = 'Se almacenará el nombre del deporte, descripción, si se juegan por equipos o de forma individual.'
$ws.Range('D72').Value = 'Mínimo'
$ws.Range('E72').Value = 'Información'
$ws.Range('F72').Value = 'Fácil'
$ws.Range('G72').Value = 'v1'
$ws.Rows(72).RowHeight = 23.85

# Row 73
$ws.Range('A73').Value = 'R72'
$ws.Range('B73').Value = 'Filtrar eventos'
$ws.Range('C73').Value = 'Se filtrarán eventos por si se practican de forma individual o por equipos'
$ws.Range('D73').Value = 'Opcional'
$ws.Range('E73').Value = 'Funcional'
$ws.Range('F73').Value = 'Media'
$ws.Range('G73').Value = 'v2'
$ws.Rows(73).RowHeight = 23.85

# Row 74
$ws.Range('A74').Value = 'R73'
$ws.Range('B74').Value = 'Buscar por fechas'
$ws.Range('C74').Value = 'Buscador para buscar eventos que se realizan desde una fecha hasta otra'
$ws.Range('D74').Value = 'Opcional'
$ws.Range('E74').Value = 'Funcional'
$ws.Range('F74').Value = 'Difícil'
$ws.Range('G74').Value = 'v3'
$ws.Rows(74).RowHeight = 23.85

# Row 75
$ws.Range('A75').Value = 'R74'
$ws.Range('B75').Value = 'Mostrar eventos del mes'
$ws.Range('C75').Value = 'En la página principal se mostrarán los eventos del mes actual'
$ws.Range('D75').Value = 'Importante'
$ws.Range('E75').Value = 'Funcional'
$ws.Range('F75').Value = 'Difícil'
$ws.Range('G75').Value = 'v2'
$ws.Rows(75).RowHeight = 23.85

# Row 76
$ws.Range('A76').Value = 'R75'
$ws.Range('B76').Value = 'Puntuación media'
$ws.Range('C76').Value = 'Cada evento muestra su puntuación media en base a las valoraciones'
$ws.Range('D76').Value = 'Importante'
$ws.Range('E76').Value = 'Funcional'
$ws.Range('F76').Value = 'Media'
$ws.Range('G76').Value = 'v2'
$ws.Rows(76).RowHeight = 23.85

# Row 77
$ws.Range('A77').Value = 'R76'
$ws.Range('B77').Value = 'Página oficial del evento'
$ws.Range('C77').Value = 'Cada evento, si dispone de ella, mostrará un enlace a su página'
$ws.Range('D77').Value = 'Importante'
$ws.Range('E77').Value = 'Funcional'
$ws.Range('F77').Value = 'Fácil'
$ws.Range('G77').Value = 'v1'
$ws.Rows(77).RowHeight = 23.85

# Row 78
$ws.Range('A78').Value = 'R77'
$ws.Range('B78').Value = 'Mostrar participantes del torneo'
$ws.Range('C78').Value = 'Cuando se ve los detalles de un evento, se mostrarán los participantes de este, de forma paginada.'
$ws.Range('D78').Value = 'Importante'
$ws.Range('E78').Value = 'Funcional'
$ws.Range('F78').Value = 'Media'
$ws.Range('G78').Value = 'v2'
$ws.Rows(78).RowHeight = 23.85

# Row 79
$ws.Range('A79').Value = 'R78'
$ws.Range('B79').Value = 'Seguir usuario'
$ws.Range('C79').Value = 'Un usuario podrá seguir a otro y ver todo los eventos que realiza'
$ws.Range('D79').Value = 'Opcional'
$ws.Range('E79').Value = 'Funcional'
$ws.Range('F79').Value = 'Difícil'
$ws.Range('G79').Value = 'v3'
$ws.Rows(79).RowHeight = 23.85

# Row 80
$ws.Range('A80').Value = 'R79'
$ws.Range('B80').Value = 'Perfil de usuario'
$ws.Range('C80').Value = 'Se mostrará una página de perfil de cada usuario'
$ws.Range('D80').Value = 'Importante'
$ws.Range('E80').Value = 'Funcional'
$ws.Range('F80').Value = 'Media'
$ws.Range('G80').Value = 'v2'
$ws.Rows(80).RowHeight = 12.8

# Row 81
$ws.Range('A81').Value = 'R80'
$ws.Range('B81').Value = 'Eventos realizados por usuario'
$ws.Range('C81').Value = 'En el perfil del usuario se podrá ver cada evento en el que se ha inscrito el usuario.'
$ws.Range('D81').Value = 'Importante'
$ws.Range('E81').Value = 'Funcional'
$ws.Range('F81').Value = 'Media'
$ws.Range('G81').Value = 'v3'
$ws.Rows(81).RowHeight = 23.85

# Row 82
$ws.Range('A82').Value = 'R81'
$ws.Range('B82').Value = 'Eventos creados por el usuario'
$ws.Range('C82').Value = 'En el perfil del usuario se podrá ver cada evento que ha creado el usuario.'
$ws.Range('D82').Value = 'Importante'
$ws.Range('E82').Value = 'Funcional'
$ws.Range('F82').Value = 'Media'
$ws.Range('G82').Value = 'v3'
$ws.Rows(82).RowHeight = 23.85

# Row 83
$ws.Range('A83').Value = 'R82'
$ws.Range('B83').Value = 'Votar comentarios'
$ws.Range('C83').Value = 'Se podrán votar los comentarios'
$ws.Range('D83').Value = 'Importante'
$ws.Range('E83').Value = 'Funcional'
$ws.Range('F83').Value = 'Fácil'
$ws.Range('G83').Value = 'v1'
$ws.Rows(83).RowHeight = 12.8

# Row 84
$ws.Range('A84').Value = 'R83'
$ws.Range('B84').Value = 'Comentar perfil de usuario'
$ws.Range('C84').Value = 'Otros usuarios podrán comentar el perfil de otro usuario'
$ws.Range('D84').Value = 'Opcional'
$ws.Range('E84').Value = 'Funcional'
$ws.Range('F84').Value = 'Media'
$ws.Range('G84').Value = 'v2'
$ws.Rows(84).RowHeight = 23.85

# Row 85
$ws.Range('A85').Value = 'R84'
$ws.Range('B85').Value = 'Compartir por redes sociales'
$ws.Range('C85').Value = 'Se podrán compartir un evento por redes sociales.'
$ws.Range('D85').Value = 'Opcional'
$ws.Range('E85').Value = 'Funcional'
$ws.Range('F85').Value = 'Difícil'
$ws.Range('G85').Value = 'v3'
$ws.Rows(85).RowHeight = 23.85

# Row 86
$ws.Range('A86').Value = 'R85'
$ws.Range('B86').Value = 'Subir cartel del evento'
$ws.Range('C86').Value = 'Al crear  un evento se podrá subir un cartel '
$ws.Range('D86').Value = 'Importante'
$ws.Range('E86').Value = 'Funcional'
$ws.Range('F86').Value = 'Media'
$ws.Range('G86').Value = 'v2'
$ws.Rows(86).RowHeight = 12.8

# Row 87
$ws.Range('A87').Value = 'R86'
$ws.Range('B87').Value = 'Formulario de contacto'
$ws.Range('C87').Value = 'La aplicación incluirá un formulario de contacto, que podrá usarse por los usuarios para realizar consultas, sugerencias o lo que deseen al administrador'
$ws.Range('D87').Value = 'Opcional'
$ws.Range('E87').Value = 'Funcional'
$ws.Range('F87').Value = 'Media'
$ws.Range('G87').Value = 'v3'
$ws.Rows(87).RowHeight = 35.05

# Row 88
$ws.Range('A88').Value = 'R87'
$ws.Range('B88').Value = 'Encuestas'
$ws.Range('C88').Value = 'En el evento se podrá añadir, si el creador lo desea, encuestas para conocer la opinión de los usuarios'
$ws.Range('D88').Value = 'Opcional'
$ws.Range('E88').Value = 'Funcional'
$ws.Range('F88').Value = 'Difícil'
$ws.Range('G88').Value = 'v3'
$ws.Rows(88).RowHeight = 23.85

# Row 89
$ws.Range('A89').Value = 'R88'
$ws.Range('B89').Value = 'Promociones'
$ws.Range('C89').Value = 'Para incentivar el pago a través de la aplicación, se mostrarán los eventos que dispongan de descuentos al pagar a través de la aplicación.'
$ws.Range('D89').Value = 'Importante'
$ws.Range('E89').Value = 'Funcional'
$ws.Range('F89').Value = 'Media'
$ws.Range('G89').Value = 'v2'
$ws.Rows(89).RowHeight = 35.05

# Row 90
$ws.Range('A90').Value = 'R89'
$ws.Range('B90').Value = 'Preguntas frecuentes'
$ws.Range('C90').Value = 'Página donde se mostrarán las preguntas frecuentes, acerca del uso de la aplicación, que pueden tener los usuarios.'
$ws.Range('D90').Value = 'Importante'
$ws.Range('E90').Value = 'Funcional'
$ws.Range('F90').Value = 'Media'
$ws.Range('G90').Value = 'v2'
$ws.Rows(90).RowHeight = 23.85

# Row 91
$ws.Range('A91').Value = 'R90'
$ws.Range('B91').Value = 'Buscar por lugar'
$ws.Range('C91').Value = 'Se mostrarán los eventos del lugar o ciudad que se indiquen.'
$ws.Range('D91').Value = 'Importante'
$ws.Range('E91').Value = 'Funcional'
$ws.Range('F91').Value = 'Media'
$ws.Range('G91').Value = 'v2'
$ws.Rows(91).RowHeight = 12.8

# Row 92
$ws.Range('A92').Value = 'R91'
$ws.Range('B92').Value = 'Videos de Youtube'
$ws.Range('C92').Value = 'En los detalles del evento, se podrán añadir videos de Youtube'
$ws.Range('D92').Value = 'Opcional'
$ws.Range('E92').Value = 'Funcional'
$ws.Range('F92').Value = 'Difícil'
$ws.Range('G92').Value = 'v3'
$ws.Rows(92).RowHeight = 12.8

# Row 93
$ws.Range('A93').Value = 'R92'
$ws.Range('B93').Value = 'Denunciar comentario'
$ws.Range('C93').Value = 'Se podrán denunciar comentarios si son irrespetuosos.'
$ws.Range('D93').Value = 'Importante'
$ws.Range('E93').Value = 'Funcional'
$ws.Range('F93').Value = 'Media'
$ws.Range('G93').Value = 'v2'
$ws.Rows(93).RowHeight = 26.85

# Row 94
$ws.Range('A94').Value = 'R93'
$ws.Range('B94').Value = 'Eliminación de comentarios con denuncias'
$ws.Range('C94').Value = 'La aplicación tendrá la capacidad de eliminar comentarios que tengan un determinado número de denuncias'
$ws.Range('D94').Value = 'Importante'
$ws.Range('E94').Value = 'Técnico'
$ws.Range('F94').Value = 'Difícil'
$ws.Range('G94').Value = 'v3'
$ws.Rows(94).RowHeight = 35.05

# Row 95
$ws.Range('A95').Value = 'R94'
$ws.Range('B95').Value = 'Artículos'
$ws.Range('C95').Value = 'Artículos informativos sobre eventos. Se almacenará, el contenido del artículo, id del usuario creador del artículo, id del evento al que se refiere.'
$ws.Range('D95').Value = 'Mínimo'
$ws.Range('E95').Value = 'Información'
$ws.Range('F95').Value = 'Fácil'
$ws.Range('G95').Value = 'v1'
$ws.Rows(95).RowHeight = 35.05

# Row 96
$ws.Range('A96').Value = 'R95'
$ws.Range('B96').Value = 'Crear artículos'
$ws.Range('C96').Value = 'Los usuarios logueados pueden crear artículos sobre los eventos que existen.'
$ws.Range('D96').Value = 'Importante'
$ws.Range('E96').Value = 'Funcional'
$ws.Range('F96').Value = 'Fácil'
$ws.Range('G96').Value = 'v1'
$ws.Rows(96).RowHeight = 23.85

# Row 97
$ws.Range('A97').Value = 'R96'
$ws.Range('B97').Value = 'Modificar artículo'
$ws.Range('C97').Value = 'El creador del artículo podrá modificarlo.'
$ws.Range('D97').Value = 'Importante'
$ws.Range('E97').Value = 'Funcional'
$ws.Range('F97').Value = 'Fácil'
$ws.Range('G97').Value = 'v1'
$ws.Rows(97).RowHeight = 12.8

# Row 98
$ws.Range('A98').Value = 'R97'
$ws.Range('B98').Value = 'Eliminar artículo'
$ws.Range('C98').Value = 'El creador del artículo podrá eliminarlo.'
$ws.Range('D98').Value = 'Importante'
$ws.Range('E98').Value = 'Funcional'
$ws.Range('F98').Value = 'Fácil'
$ws.Range('G98').Value = 'v1'
$ws.Rows(98).RowHeight = 12.8

# Row 99
$ws.Range('A99').Value = 'R98'
$ws.Range('B99').Value = 'Denunciar artículo'
$ws.Range('C99').Value = 'Los usuarios podrán denunciar un artículo'
$ws.Range('D99').Value = 'Importante'
$ws.Range('E99').Value = 'Funcional'
$ws.Range('F99').Value = 'Media'
$ws.Range('G99').Value = 'v2'
$ws.Rows(99).RowHeight = 12.8

# Row 100
$ws.Range('A100').Value = 'R99'
$ws.Range('B100').Value = 'Eliminación de artículos con denuncias'
$ws.Range('C100').Value = 'La aplicación tendrá la capacidad de eliminar artículos que tengan un determinado número de denuncias'
$ws.Range('D100').Value = 'Importante'
$ws.Range('E100').Value = 'Técnico'
$ws.Range('F100').Value = 'Difícil'
$ws.Range('G100').Value = 'v3'
$ws.Rows(100).RowHeight = 35.05

# Row 101
$ws.Range('A101').Value = 'R100'
$ws.Range('B101').Value = 'Validación de usuarios'
$ws.Range('C101').Value = 'Los usuarios recibirán un email para validar su cuenta.'
$ws.Range('D101').Value = 'Importante'
$ws.Range('E101').Value = 'Funcional'
$ws.Range('F101').Value = 'Media'
$ws.Range('G101').Value = 'v1'
$ws.Rows(101).RowHeight = 12.8

# Row 102
$ws.Range('A102').Value = 'R101'
$ws.Range('B102').Value = 'Restricción del uso para usuarios no validados'
$ws.Range('C102').Value = 'Los usuarios no validados no podrán acceder a la aplicación.'
$ws.Range('D102').Value = 'Importante'
$ws.Range('E102').Value = 'Funcional'
$ws.Range('F102').Value = 'Media'
$ws.Range('G102').Value = 'v1'
$ws.Rows(102).RowHeight = 35.05

# Row 103
$ws.Range('A103').Value = 'R102'
$ws.Range('B103').Value = 'Eliminación de usuarios no validados'
$ws.Range('C103').Value = 'Pasado un determinado tiempo, si el usuario no se ha validado, la aplicación borrará su cuenta.'
$ws.Range('D103').Value = 'Importante'
$ws.Range('E103').Value = 'Funcional'
$ws.Range('F103').Value = 'Difícil'
$ws.Range('G103').Value = 'v3'
$ws.Rows(103).RowHeight = 23.85

# Row 104
$ws.Range('A104').Value = 'R103'
$ws.Range('B104').Value = 'Mapas'
$ws.Range('C104').Value = 'En el evento se mostrará la ubicación de este mediante Google Maps'
$ws.Range('D104').Value = 'Opcional'
$ws.Range('E104').Value = 'Funcional'
$ws.Range('F104').Value = 'Difícil'
$ws.Range('G104').Value = 'v3'
$ws.Rows(104).RowHeight = 23.85

# 4. Update selection to match final state
$ws.Range('A105').Select()

Write-Host "done"